$d = $word.ActiveDocument

function Replace-Exact([string]$findText, [string]$replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

function SetText-NoAutocorrect([string]$findText, [string]$newText) {
    # Finds findText and overwrites the matched Range's .Text directly, which
    # (unlike Find.Execute's Replace) does not trigger AutoCorrect / smart-quote
    # substitution, so straight apostrophes etc. survive untouched.
    $rng = $d.Content
    $rng.Find.Execute($findText) | Out-Null
    $rng.Text = $newText
}

# 1. "RevoScaleR" stays bold/own run; merge the rest of the bullet into one run.
Replace-Exact " reads data from a variety of sources + can prep data, run descriptive statistics + statistical tests" " reads data from a variety of sources + can prep data, run descriptive statistics + statistical tests"

# 2. "Can take a sample..." bullet -> merge entirely into one run.
#    (Contains a straight apostrophe in "R's" - merge using a placeholder then
#    restore the apostrophe without triggering AutoCorrect.)
Replace-Exact "Can take a sample from our data + w/ that sample, have a data frame that can be used by all of R's packages for performing various analysis." "Can take a sample from our data + w/ that sample, have a data frame that can be used by all of R@APOS@s packages for performing various analysis."
SetText-NoAutocorrect "R@APOS@s packages" "R's packages"

# 3. "Bread + butter of RevoScaleR = the " merges (before the bold "analytics algorithms").
Replace-Exact "Bread + butter of RevoScaleR = the " "Bread + butter of RevoScaleR = the "

# 3b. " = regression + classification ... random forests, +" merges (after bold "analytics algorithms").
Replace-Exact " = regression + classification algorithms such as linear models, logistic regression, decision trees, ensemble models such as random forests, +" " = regression + classification algorithms such as linear models, logistic regression, decision trees, ensemble models such as random forests, +"

# 3c. " k-means algorithms" merges separately (stays its own run, not merged w/ neighbours).
Replace-Exact " k-means algorithms" " k-means algorithms"

# 4. "All these algorithms have counterparts ... they're " merges (before bold "parallel").
Replace-Exact "All these algorithms have counterparts in open source R functions but strength of MSR algorithms = they’re " "All these algorithms have counterparts in open source R functions but strength of MSR algorithms = they’re "

# 4b. " = makes it so we can run algorithms on very large data sets in a scalable fashion" merges (after bold "parallel").
Replace-Exact " = makes it so we can run algorithms on very large data sets in a scalable fashion" " = makes it so we can run algorithms on very large data sets in a scalable fashion"

# 5. "In addition to being parallel + scalable, these algorithms can also run inside of PROD environments" merges entirely.
Replace-Exact "In addition to being parallel + scalable, these algorithms can also run inside of PROD environments" "In addition to being parallel + scalable, these algorithms can also run inside of PROD environments"

# 6. "Ex: What's h" merges (keeps lastRenderedPageBreak on first run). No apostrophe in search/replace text needed.
Replace-Exact "Ex: What" "Ex: What"

# 7. "Left = Code to point to some data on a local Linux or Windows machine" - merge last two runs.
Replace-Exact "to some data on a local Linux or Windows machine" "to some data on a local Linux or Windows machine"

# 8. "Bottom = Take..." -> "T"+"ak" merges to "Tak"; " data "+"+" merges to " data +".
Replace-Exact "Tak" "Tak"
Replace-Exact " data +" " data +"

# 9. "3 main benefits of RevoScaleR package." merges (underline run).
Replace-Exact "3 main benefits of RevoScaleR package." "3 main benefits of RevoScaleR package."

# 10. "1) Even w/ ... counterparts " merges entirely.
#     (Contains a straight apostrophe in "RevoScaleR's" - use placeholder trick.)
Replace-Exact "1) Even w/ enough memory to load data as a data frame into an R session, can use RevoScaleR's parallel algorithms to run analytics on that data much faster than w/ open source counterparts " "1) Even w/ enough memory to load data as a data frame into an R session, can use RevoScaleR@APOS@s parallel algorithms to run analytics on that data much faster than w/ open source counterparts "
SetText-NoAutocorrect "RevoScaleR@APOS@s parallel" "RevoScaleR's parallel"

# 11. "2) If data is too large ... disks." merges entirely.
Replace-Exact "2) If data is too large to fit in available memory, can still use RevoScaleR algorithms just as before by simply pointing to data sitting on disks." "2) If data is too large to fit in available memory, can still use RevoScaleR algorithms just as before by simply pointing to data sitting on disks."

# 12. "RevoScaleR operates by loading ... changed)" merges entirely.
Replace-Exact "RevoScaleR operates by loading data into R session as a data frame but only a chunk at a time (default = 500K rows but can be changed)" "RevoScaleR operates by loading data into R session as a data frame but only a chunk at a time (default = 500K rows but can be changed)"

# 13. "By doing so RevoScaleR ... handled." merges entirely.
Replace-Exact "By doing so RevoScaleR can simply load data 1 chunk at a time, process it, move on to the next chunk, + keep doing this until all data has been handled." "By doing so RevoScaleR can simply load data 1 chunk at a time, process it, move on to the next chunk, + keep doing this until all data has been handled."

# 14. "The fact that its algorithms are parallel is what makes this possible." merges entirely.
Replace-Exact "The fact that its algorithms are parallel is what makes this possible." "The fact that its algorithms are parallel is what makes this possible."

# 15. "3) Can take code ... code structure." merges entirely.
Replace-Exact "3) Can take code + deploy it inside the PROD environment (Hadoop cluster or SQL Server database) w/ very little changes made to the code structure." "3) Can take code + deploy it inside the PROD environment (Hadoop cluster or SQL Server database) w/ very little changes made to the code structure."

# 16. "Microsoft R Client = lightweight version ... single machi" merges (keeps "ne" separate, and bold lead-in separate).
Replace-Exact " lightweight version of MRS (not meant to be use as a PROD environment for the MRS on a single machi" " lightweight version of MRS (not meant to be use as a PROD environment for the MRS on a single machi"

# 18. "to C:Program Files/" is re-split into "to " + "C:Program" (proofErr gramStart/End) + " Files/" -
#     identical visible text, so no text change is needed there via COM.
#     Actual content change: "/ R_SERVER (default directory)" -> "/R_SERVER (default directory)" (drop the space).
Replace-Exact "/ R_SERVER (default directory)" "/R_SERVER (default directory)"
